# Apply crypto price/volume updates scraped by the GitHub Actions job.
# (Wed Sep 20 19:58:36 UTC 2023 "Updated cryptos list" commit.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.995.33"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.619.90"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.98"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0624"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.17"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "1.640.61"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.12"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.541"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.45"
$ws.Range("E15").Value = "  -4.74%  "
$ws.Range("D16").Value = "26.989.26"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "0.0₃0741"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.83"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("E22").Value = "  -5.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.99"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.51"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.27"
$ws.Range("E26").Value = "  -3.93%  "
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.55"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.35"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("D33").Value = "1.332.38"
$ws.Range("E33").Value = "  +5.72%  "
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.543"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.845"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.801"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.23"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.55"
$ws.Range("E42").Value = "  +4.37%  "
$ws.Range("E43").Value = "  -3.58%  "
$ws.Range("D44").Value = "1.757.14"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.32"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("E48").Value = "  +22.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0982"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("E51").Value = "  -1.09%  "
